# Regenerate the "K" column (column G) values for Sheet1.
# The former "Strike#" values in column G are replaced with actual
# strikeout counts (K) pulled from the regenerated save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number (data rows 2..68) -> new K value.
$newK = [ordered]@{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 3
    7  = 2
    8  = 1
    9  = 2
    10 = 3
    11 = 0
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 1
    29 = 2
    30 = 0
    31 = 0
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 2
    44 = 3
    45 = 0
    46 = 1
    47 = 1
    48 = 2
    49 = 0
    50 = 2
    51 = 0
    52 = 1
    53 = 0
    54 = 0
    55 = 2
    56 = 0
    57 = 1
    58 = 1
    59 = 2
    60 = 2
    61 = 2
    62 = 2
    63 = 1
    64 = 3
    65 = 2
    66 = 2
    67 = 2
    68 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
